$p = $ppt.ActivePresentation

# --- Slide 1: merge "18 " + "Sep 2019" runs into a single "18 Sep 2019" run ---
$s1 = $p.Slides.Item(1)
$sh1 = $s1.Shapes.Item(2)
$tr1 = $sh1.TextFrame.TextRange
$dateParagraph = $tr1.Paragraphs(2)
# Force a real content change first so the engine collapses the two runs into
# one, then restore the desired final text on that single run. The
# placeholder intentionally shares no characters with the final text so the
# text-diff engine cannot "helpfully" keep a multi-run split.
$dateParagraph.Text = "XXXXXXXXXXXXXXXXXXXXXXXXXXX"
$dateParagraph.Text = "18 Sep 2019"

# --- Slide 6: split "Not so large (~5K)" into "Not so large (~" + "50M)" ---
$s6 = $p.Slides.Item(6)
$sh6 = $s6.Shapes.Item(2)
$tr6 = $sh6.TextFrame.TextRange
$sizeParagraph = $tr6.Paragraphs(2)
# "Not so large (~" is the first 15 characters; the remaining "5K)" becomes "50M)".
$tail = $sizeParagraph.Characters(16, 3)
$tail.Text = "50M)"
